$d = $word.ActiveDocument

# --- Step 1: Remove the "Meta description: ..." paragraph (paragraph 2) ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# --- Step 2: Insert a new bold paragraph "Play Cleopatra II Slot for Free - Maximum Payout 10,000x"
#             right before the final ("Create a feature image...") paragraph. ---
$n = $d.Paragraphs.Count
$pPrev = $d.Paragraphs.Item($n - 1)
$endR = $d.Range($pPrev.Range.End, $pPrev.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cleopatra II Slot for Free - Maximum Payout 10,000x</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endR.InsertXML($xml)

# Split the paragraph right before the (still-present) old italic text so the inserted
# bold text becomes its own paragraph and the trailing text keeps its own leading empty run.
$splitR = $d.Content
$splitR.Find.Execute("Create a feature image", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$collapsedR = $d.Range($splitR.Start, $splitR.Start)
$collapsedR.InsertParagraphBefore() | Out-Null

# --- Step 3: Replace the text of the final paragraph (keeping its italic formatting) ---
$old = "Create a feature image with a cartoon-style happy Maya warrior with glasses to fit the game Cleopatra II. The image should have a fun and vibrant design, showcasing the warrior in a lively and colorful way. The Maya warrior should be wearing a traditional outfit and headdress, holding a treasure chest or other valuable item related to the game. The glasses add a playful touch to the image, adding a modern and unique element to the traditional setting. The image should be eye-catching and engaging, drawing players in and making them excited to play the game."
$new = "Read our review of Cleopatra II slot game. Play for free and win up to 10,000 times your bet with maximum 180 free spins."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
